# Edit script: add the "24. 8. 2021" wave of data to both worksheets
# (sheet "data" -> new column AH, sheet "pocetR" -> new column AG),
# and refresh the "aktualizace" date in each footer title row.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "data"  (percentages) -- new column AH (col 34)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Header cell AH1 needs the same style as the other date headers (AG1),
# so copy formatting across before setting its text.
$ws1.Range("AG1").Copy()
$ws1.Range("AH1").PasteSpecial(-4122)
$ws1.Range("AH1").Value = "24. 8. 2021"

$data1 = @(0.64, 0.06, 0.49, 0.09, 0.6, 0.07000000000000001, 0.76, 0.04, 0.59, 0.07000000000000001, 0.63, 0.07000000000000001, 0.77, 0.04, 0.49, 0.07000000000000001, 0.53, 0.09, 0.73, 0.05, 0.68, 0.06, 0.86, 0.02, 0.87, 0.03, 0.41, 0.13, 0.65, 0.07000000000000001, 0.6899999999999999, 0.06, 0.72, 0.04, 0.68, 0.06, 0.59, 0.07000000000000001, 0.58, 0.06, 0.65, 0.06, 0.8, 0.05, 0.39, 0.11, 0.55, 0.08, 0.73, 0.04, 0.61, 0.07000000000000001, 0.62, 0.06, 0.63, 0.05, 0.76, 0.04, 0.64, 0.05, 0.71, 0.03, 0.68, 0.05, 0.8, 0.05, 0.55, 0.08, 0.6899999999999999, 0.07000000000000001, 0.6899999999999999, 0.06, 0.61, 0.06)

$r = 2
foreach ($v in $data1) {
    $ws1.Cells.Item($r, 34).Value = $v
    $r = $r + 1
}

# Footer title row (76): bump the "aktualizace" date
$ws1.Range("A76").Value = "Život během pandemie, Imunizace, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# ------------------------------------------------------------------
# Sheet 2: "pocetR"  (sample sizes) -- new column AG (col 33)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AG1: copy style from AF1, then set the date label
$ws2.Range("AF1").Copy()
$ws2.Range("AG1").PasteSpecial(-4122)
$ws2.Range("AG1").Value = "24. 8. 2021"

$data2 = @(1901, 452, 709, 740, 872, 654, 375, 445, 458, 258, 427, 196, 117, 200, 729, 583, 263, 937, 964, 235, 370, 332, 217, 339, 408, 988, 442, 220, 251, 263, 257, 354, 249, 515, 314, 350, 1237)

$r = 2
foreach ($v in $data2) {
    $ws2.Cells.Item($r, 33).Value = $v
    $r = $r + 1
}

# Footer title row (39): bump the "aktualizace" date
$ws2.Range("A39").Value = "Život během pandemie, Imunizace, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"

# AG39 stays an empty placeholder cell, matching the blank cells already
# present across the rest of that footer row (B39:AF39).
$ws2.Cells.Item(39, 33).Value = ""

